$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the calibration-factor formula in row 14: use the C14 cell
#     (number of divisions) instead of the hard-coded literal 12, and
#     correct the actual measured value to 14 divisions. ---
$ws.Range("C14").Value = 14
$ws.Range("F14").Formula = "=546/C14"

# --- Finish the lab: compute, for every dark/light ring pair (rows 4-10),
#     the radius-squared values and the "delta" used for the wavelength
#     calculation (formula 2.2), laid out in the new block B17:N23. ---

# Row 17 holds the first data point plus the summary formulas (J17, M17, N17)
$ws.Range("B17").Formula = '=B4+C4/100'
$ws.Range("C17").Formula = '=- (B17-3.61) * 0.1'
$ws.Range("D17").Formula = '=C17*C17'
$ws.Range("E17").Formula = '=E4+F4/100'
$ws.Range("F17").Formula = '=- (E17-3.61) * 0.1'
$ws.Range("G17").Formula = '=F17*F17'
$ws.Range("H17").Formula = '=G17-D18'
$ws.Range("J17").Formula = '=AVERAGE(H17:H22)'
$ws.Range("K17").Formula = '=-H17-$J$17'
$ws.Range("L17").Formula = '=K17*K17'
$ws.Range("M17").Formula = '=SUM(L17:L23)'
$ws.Range("N17").Formula = '=SQRT(M17)'

# Rows 18-23 repeat the same per-row formulas (B..G, K, L), each sourced
# from the corresponding measurement row 5-10; H is only needed through
# row 22 (it references the next row's D).
for ($row = 18; $row -le 23; $row++) {
    $src = $row - 13
    $ws.Range("B$row").Formula = "=B$src+C$src/100"
    $ws.Range("C$row").Formula = "=- (B$row-3.61) * 0.1"
    $ws.Range("D$row").Formula = "=C$row*C$row"
    $ws.Range("E$row").Formula = "=E$src+F$src/100"
    $ws.Range("F$row").Formula = "=- (E$row-3.61) * 0.1"
    $ws.Range("G$row").Formula = "=F$row*F$row"
    if ($row -le 22) {
        $next = $row + 1
        $ws.Range("H$row").Formula = "=G$row-D$next"
    }
    $ws.Range("K$row").Formula = '=-H' + $row + '-$J$17'
    $ws.Range("L$row").Formula = "=K$row*K$row"
}

# --- Restore the view: scrolled so row 4 is at the top, with the new
#     C24 cell (just below the finished table) selected. ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C24").Select()
